$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 102.36364
$ws.Range("I11").Value = 102.36364
$ws.Range("K11").Value = 102.36364
$ws.Range("M11").Value = 37.63636
$ws.Range("H21").Value = 6500
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 6500
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H26").Value = 25000
$ws.Range("J26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("N26").Value = -25688
$ws.Range("H32").Value = 4856.2856
$ws.Range("J32").Value = 4832.6665
$ws.Range("L32").Value = 4832.6665
$ws.Range("N32").Value = -5484.6665
$ws.Range("H40").Value = 2119.2
$ws.Range("J40").Value = 3500
$ws.Range("L40").Value = 3500
$ws.Range("N40").Value = -3850
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H80").Value = 901.7368
$ws.Range("I80").Value = 318.57144
$ws.Range("J80").Value = 1241.9166
$ws.Range("K80").Value = 955.71432
$ws.Range("L80").Value = 3725.7498
$ws.Range("M80").Value = 42.28567999999996
$ws.Range("N80").Value = -5721.7498
$ws.Range("H83").Value = 901.7368
$ws.Range("I83").Value = 318.57144
$ws.Range("J83").Value = 1241.9166
$ws.Range("K83").Value = 2867.14296
$ws.Range("L83").Value = 11177.2494
$ws.Range("M83").Value = 2124.85704
$ws.Range("N83").Value = -21161.2494
$ws.Range("H86").Value = 3225.2
$ws.Range("I86").Value = 3002
$ws.Range("K86").Value = 3002
$ws.Range("M86").Value = -1879
$ws.Range("H89").Value = 3225.2
$ws.Range("I89").Value = 3002
$ws.Range("K89").Value = 15010
$ws.Range("M89").Value = -9394
$ws.Range("H98").Value = 8820.477000000001
$ws.Range("I98").Value = 10029.111
$ws.Range("K98").Value = 10029.111
$ws.Range("M98").Value = -8531.111000000001
$ws.Range("H100").Value = 5001.1665
$ws.Range("I100").Value = 2165
$ws.Range("K100").Value = 2165
$ws.Range("M100").Value = -1624
$ws.Range("H112").Value = 1771.6111
$ws.Range("I112").Value = 1298
$ws.Range("J112").Value = 1953.7693
$ws.Range("K112").Value = 3894
$ws.Range("L112").Value = 5861.3079
$ws.Range("M112").Value = -2786
$ws.Range("N112").Value = -8077.3079
$ws.Range("H122").Value = 8820.477000000001
$ws.Range("I122").Value = 10029.111
$ws.Range("K122").Value = 30087.333
$ws.Range("M122").Value = -27637.333
$ws.Range("H137").Value = 1682.2
$ws.Range("I137").Value = 1286.4445
$ws.Range("J137").Value = 2699.8572
$ws.Range("K137").Value = 3859.3335
$ws.Range("L137").Value = 8099.571599999999
$ws.Range("M137").Value = -1309.3335
$ws.Range("N137").Value = -13199.5716
$ws.Range("H138").Value = 6071.4585
$ws.Range("I138").Value = 8106.643
$ws.Range("J138").Value = 5580.207
$ws.Range("K138").Value = 24319.929
$ws.Range("L138").Value = 16740.621
$ws.Range("M138").Value = -19179.929
$ws.Range("N138").Value = -27020.621
$ws.Range("H141").Value = 6558.963
$ws.Range("I141").Value = 2924.5
$ws.Range("K141").Value = 8773.5
$ws.Range("M141").Value = -3593.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 545
$ws.Range("I21").Value = 393.33334
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 393.33334
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = -19.33334000000002
$ws.Range("N21").Value = -1748
$ws.Range("H23").Value = 29761.904
$ws.Range("J23").Value = 25000
$ws.Range("L23").Value = 25000
$ws.Range("N23").Value = -25518
$ws.Range("H43").Value = 12087
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12087
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12087
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -12713
$ws.Range("H45").Value = 9066.538
$ws.Range("J45").Value = 1610.75
$ws.Range("L45").Value = 1610.75
$ws.Range("N45").Value = -2364.75
$ws.Range("H61").Value = 5565.6875
$ws.Range("I61").Value = 5669.1333
$ws.Range("J61").Value = 4014
$ws.Range("K61").Value = 5669.1333
$ws.Range("L61").Value = 4014
$ws.Range("M61").Value = -5457.1333
$ws.Range("N61").Value = -4438
$ws.Range("H63").Value = 6825
$ws.Range("J63").Value = 9900
$ws.Range("L63").Value = 9900
$ws.Range("N63").Value = -11272
$ws.Range("H66").Value = 6825
$ws.Range("J66").Value = 9900
$ws.Range("L66").Value = 49500
$ws.Range("N66").Value = -56364
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51748
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -158736
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H97").Value = 627.2222
$ws.Range("I97").Value = 520.7143
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 520.7143
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -24.71429999999998
$ws.Range("N97").Value = -1992
$ws.Range("H132").Value = 1832.9131
$ws.Range("I132").Value = 1802.591
$ws.Range("K132").Value = 5407.772999999999
$ws.Range("M132").Value = -2877.772999999999
$ws.Range("H136").Value = 5565.6875
$ws.Range("I136").Value = 5669.1333
$ws.Range("J136").Value = 4014
$ws.Range("K136").Value = 17007.3999
$ws.Range("L136").Value = 12042
$ws.Range("M136").Value = -14457.3999
$ws.Range("N136").Value = -17142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7746.857
$ws.Range("I20").Value = 11681.556
$ws.Range("J20").Value = 664.4
$ws.Range("K20").Value = 11681.556
$ws.Range("L20").Value = 664.4
$ws.Range("M20").Value = -11434.556
$ws.Range("N20").Value = -1158.4
$ws.Range("H63").Value = 50001
$ws.Range("J63").Value = 50001
$ws.Range("L63").Value = 50001
$ws.Range("N63").Value = -51373
$ws.Range("H66").Value = 50001
$ws.Range("J66").Value = 50001
$ws.Range("L66").Value = 150003
$ws.Range("N66").Value = -156867
$ws.Range("H75").Value = 20000
$ws.Range("I75").Value = 20000
$ws.Range("K75").Value = 20000
$ws.Range("M75").Value = -19064
$ws.Range("H78").Value = 20000
$ws.Range("I78").Value = 20000
$ws.Range("K78").Value = 60000
$ws.Range("M78").Value = -55320
$ws.Range("H86").Value = 3186242.2
$ws.Range("I86").Value = 36666
$ws.Range("J86").Value = 6335818.5
$ws.Range("K86").Value = 36666
$ws.Range("L86").Value = 6335818.5
$ws.Range("M86").Value = -35543
$ws.Range("N86").Value = -6338064.5
$ws.Range("H89").Value = 3186242.2
$ws.Range("I89").Value = 36666
$ws.Range("J89").Value = 6335818.5
$ws.Range("K89").Value = 183330
$ws.Range("L89").Value = 31679092.5
$ws.Range("M89").Value = -177714
$ws.Range("N89").Value = -31690324.5
$ws.Range("H94").Value = 2272.0476
$ws.Range("I94").Value = 1751
$ws.Range("K94").Value = 1751
$ws.Range("M94").Value = -1300
$ws.Range("H99").Value = 4376.067
$ws.Range("I99").Value = 4059.182
$ws.Range("J99").Value = 5247.5
$ws.Range("K99").Value = 4059.182
$ws.Range("L99").Value = 5247.5
$ws.Range("M99").Value = -2561.182
$ws.Range("N99").Value = -8243.5
$ws.Range("H105").Value = 3885.1064
$ws.Range("I105").Value = 4467.5947
$ws.Range("J105").Value = 1729.9
$ws.Range("K105").Value = 4467.5947
$ws.Range("L105").Value = 1729.9
$ws.Range("M105").Value = -2720.5947
$ws.Range("N105").Value = -5223.9
$ws.Range("H134").Value = 2131.077
$ws.Range("I134").Value = 2089
$ws.Range("J134").Value = 2271.3333
$ws.Range("K134").Value = 6267
$ws.Range("L134").Value = 6813.999899999999
$ws.Range("M134").Value = -3732
$ws.Range("N134").Value = -11883.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1386.3334
$ws.Range("J22").Value = 1878.9
$ws.Range("L22").Value = 1878.9
$ws.Range("N22").Value = -2578.9
$ws.Range("H42").Value = 9500
$ws.Range("I42").Value = 9500
$ws.Range("K42").Value = 9500
$ws.Range("M42").Value = -8907
$ws.Range("H58").Value = 2331.5715
$ws.Range("I58").Value = 2304.4
$ws.Range("J58").Value = 2399.5
$ws.Range("K58").Value = 2304.4
$ws.Range("L58").Value = 2399.5
$ws.Range("M58").Value = -2101.4
$ws.Range("N58").Value = -2805.5
$ws.Range("H68").Value = 64999.668
$ws.Range("J68").Value = 67500
$ws.Range("L68").Value = 67500
$ws.Range("N68").Value = -68998
$ws.Range("H71").Value = 64999.668
$ws.Range("J71").Value = 67500
$ws.Range("L71").Value = 202500
$ws.Range("N71").Value = -209988
$ws.Range("H74").Value = 34500
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51748
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 34500
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -158736
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H81").Value = 75000
$ws.Range("J81").Value = 75000
$ws.Range("L81").Value = 75000
$ws.Range("N81").Value = -76996
$ws.Range("H84").Value = 75000
$ws.Range("J84").Value = 75000
$ws.Range("L84").Value = 225000
$ws.Range("N84").Value = -234984
$ws.Range("H88").Value = 1000
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 1000
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 26722.75
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H112").Value = 61875
$ws.Range("J112").Value = 61875
$ws.Range("L112").Value = 61875
$ws.Range("N112").Value = -64829
$ws.Range("H132").Value = 1908.3334
$ws.Range("I132").Value = 1885.9524
$ws.Range("K132").Value = 5657.857199999999
$ws.Range("M132").Value = -3127.857199999999
$ws.Range("H134").Value = 4559.5
$ws.Range("I134").Value = 3955
$ws.Range("K134").Value = 11865
$ws.Range("M134").Value = -9330
$ws.Range("H136").Value = 2331.5715
$ws.Range("I136").Value = 2304.4
$ws.Range("J136").Value = 2399.5
$ws.Range("K136").Value = 6913.200000000001
$ws.Range("L136").Value = 7198.5
$ws.Range("M136").Value = -4363.200000000001
$ws.Range("N136").Value = -12298.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 65.166664
$ws.Range("I2").Value = 49.75
$ws.Range("K2").Value = 298.5
$ws.Range("M2").Value = -185.5
$ws.Range("H112").Value = 6912.5
$ws.Range("I112").Value = 1825
$ws.Range("J112").Value = 12000
$ws.Range("K112").Value = 5475
$ws.Range("L112").Value = 36000
$ws.Range("M112").Value = -4367
$ws.Range("N112").Value = -38216
$ws.Range("H121").Value = 2599.4375
$ws.Range("I121").Value = 512.6667
$ws.Range("J121").Value = 5282.4287
$ws.Range("K121").Value = 1538.0001
$ws.Range("L121").Value = 15847.2861
$ws.Range("M121").Value = -228.0001
$ws.Range("N121").Value = -18467.2861
$ws.Range("H140").Value = 1294
$ws.Range("I140").Value = 1294
$ws.Range("K140").Value = 3882
$ws.Range("M140").Value = 1298

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 35000
$ws.Range("J26").Value = 35000
$ws.Range("L26").Value = 35000
$ws.Range("N26").Value = -35560
$ws.Range("H43").Value = 7000
$ws.Range("I43").Value = 18000
$ws.Range("J43").Value = 6153.846
$ws.Range("K43").Value = 18000
$ws.Range("L43").Value = 6153.846
$ws.Range("M43").Value = -17849
$ws.Range("N43").Value = -6455.846
$ws.Range("H50").Value = 35000
$ws.Range("J50").Value = 35000
$ws.Range("L50").Value = 35000
$ws.Range("N50").Value = -35996
$ws.Range("H62").Value = 38000
$ws.Range("I62").Value = 38000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 38000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -37314
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 38000
$ws.Range("I65").Value = 38000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 114000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -110568
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 4712.5
$ws.Range("J70").Value = 4496.7
$ws.Range("L70").Value = 4496.7
$ws.Range("N70").Value = -5036.7
$ws.Range("H73").Value = 4712.5
$ws.Range("J73").Value = 4496.7
$ws.Range("L73").Value = 4496.7
$ws.Range("N73").Value = -6368.7
$ws.Range("H74").Value = 50001
$ws.Range("J74").Value = 50001
$ws.Range("L74").Value = 50001
$ws.Range("N74").Value = -51873
$ws.Range("H77").Value = 50001
$ws.Range("J77").Value = 50001
$ws.Range("L77").Value = 150003
$ws.Range("N77").Value = -159363
$ws.Range("H80").Value = 5621.9
$ws.Range("I80").Value = 3243.8
$ws.Range("J80").Value = 8000
$ws.Range("K80").Value = 3243.8
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = -2245.8
$ws.Range("N80").Value = -9996
$ws.Range("H83").Value = 5621.9
$ws.Range("I83").Value = 3243.8
$ws.Range("J83").Value = 8000
$ws.Range("K83").Value = 16219
$ws.Range("L83").Value = 40000
$ws.Range("M83").Value = -11227
$ws.Range("N83").Value = -49984
$ws.Range("H97").Value = 517.6129
$ws.Range("I97").Value = 374.96155
$ws.Range("K97").Value = 374.96155
$ws.Range("M97").Value = 121.03845
$ws.Range("H132").Value = 3190.077
$ws.Range("I132").Value = 3190.077
$ws.Range("K132").Value = 9570.231
$ws.Range("M132").Value = -7040.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2751.647
$ws.Range("I22").Value = 2330.6667
$ws.Range("J22").Value = 2981.2727
$ws.Range("K22").Value = 2330.6667
$ws.Range("L22").Value = 2981.2727
$ws.Range("M22").Value = -2035.6667
$ws.Range("N22").Value = -3571.2727
$ws.Range("H27").Value = 2751.647
$ws.Range("I27").Value = 2330.6667
$ws.Range("J27").Value = 2981.2727
$ws.Range("K27").Value = 2330.6667
$ws.Range("L27").Value = 2981.2727
$ws.Range("M27").Value = -2223.6667
$ws.Range("N27").Value = -3195.2727
$ws.Range("H46").Value = 3438.3333
$ws.Range("I46").Value = 2333.3333
$ws.Range("J46").Value = 3990.8333
$ws.Range("K46").Value = 2333.3333
$ws.Range("L46").Value = 3990.8333
$ws.Range("M46").Value = -2145.3333
$ws.Range("N46").Value = -4366.8333
$ws.Range("H55").Value = 2191.6
$ws.Range("I55").Value = 986.6667
$ws.Range("K55").Value = 986.6667
$ws.Range("M55").Value = -813.6667
$ws.Range("H61").Value = 3610.7407
$ws.Range("I61").Value = 2527.2778
$ws.Range("K61").Value = 2527.2778
$ws.Range("M61").Value = -2325.2778
$ws.Range("H82").Value = 3729.6
$ws.Range("I82").Value = 1298.6
$ws.Range("K82").Value = 1298.6
$ws.Range("M82").Value = -937.5999999999999
$ws.Range("H85").Value = 3729.6
$ws.Range("I85").Value = 1298.6
$ws.Range("K85").Value = 1298.6
$ws.Range("M85").Value = -50.59999999999991
$ws.Range("H93").Value = 4067.2
$ws.Range("I93").Value = 1334.5
$ws.Range("K93").Value = 1334.5
$ws.Range("M93").Value = -86.5
$ws.Range("H100").Value = 6966.4165
$ws.Range("I100").Value = 4799.5713
$ws.Range("K100").Value = 4799.5713
$ws.Range("M100").Value = -4258.5713
$ws.Range("H113").Value = 3610.7407
$ws.Range("I113").Value = 2527.2778
$ws.Range("K113").Value = 2527.2778
$ws.Range("M113").Value = -357.2777999999998
$ws.Range("H128").Value = 59500
$ws.Range("J128").Value = 59500
$ws.Range("L128").Value = 59500
$ws.Range("N128").Value = -69460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H81").Value = 1784.8235
$ws.Range("I81").Value = 1289.4667
$ws.Range("K81").Value = 2578.9334
$ws.Range("M81").Value = -1517.9334
$ws.Range("H84").Value = 1784.8235
$ws.Range("I84").Value = 1289.4667
$ws.Range("K84").Value = 12894.667
$ws.Range("M84").Value = -7590.666999999999
$ws.Range("H96").Value = 8000
$ws.Range("I96").Value = 8000
$ws.Range("J96").Value = 8000
$ws.Range("K96").Value = 8000
$ws.Range("L96").Value = 8000
$ws.Range("M96").Value = -6627
$ws.Range("N96").Value = -10746
$ws.Range("H97").Value = 18230
$ws.Range("J97").Value = 18230
$ws.Range("L97").Value = 18230
$ws.Range("N97").Value = -20212
